$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.364.05'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.54%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.874.56'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.77%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7122'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.28%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.93'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.67%  '

# Row 7
$ws.Range("E7").Value = '  -0.05%  '

# Row 8
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07796'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.09%  '

# Row 9
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3111'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.24%  '

# Row 10
$ws.Range("E10").Value = '  +2.05%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08440'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.81%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.875.85'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.79%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.241'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.94%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7136'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.50%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.17'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.02%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.369.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.57%  '

# Row 17
$ws.Range("E17").Value = '  +2.28%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008234'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.44%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.94'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.54%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.98%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.122.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.08%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.04%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.780'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.80%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.03%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1598'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.54%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.35%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.070'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.10%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.58%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.511'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.06%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.422'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.60%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.327'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.62%  '

# Row 32
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.289'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.30%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05306'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.69%  '

# Row 34
$ws.Range("E34").Value = '  +0.89%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.180'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.38%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7442'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.12%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.698'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.50%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01868'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.39%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.223.31'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.08%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.733'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.49%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.519'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.36%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '110.78'
$ws.Range("D42").Style = "Normal"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8919'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.16%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.55%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.01%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.019.91'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.21%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.814'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.82%  '

# Row 48
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000123'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.87%  '

# Row 49
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5215'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.67%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.435'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.31%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4323'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.53%  '
